# Regenerate merged AHB files
# 1. Rename header columns: *_old -> *_FV2410, *_new -> *_FV2504
# 2. Convert the data range into an Excel Table ("Table1")
# 3. Freeze the header row (top row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rename header row cells ------------------------------------
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace "_old$", "_FV2410"
        $newVal = $newVal -replace "_new$", "_FV2504"
        $cell.Value2 = $newVal
    }
}

# --- Step 2: create the table over the used range ------------------------
$tableRange = $ws.Range("A1:U60")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# --- Step 3: freeze the top row -------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
